$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.041.77"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.496.55"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "2.520.32"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "2.941.35"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "58.887.09"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.513.98"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
